$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column-A cells beyond the original range (rows 24-25) inherit the bold/bordered
# style used throughout column A (same style as the header row / A2:A23).
$ws.Range("A2").Copy()
$ws.Range("A24:A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write the full, final table (row 2 is a newly inserted record; row 11 is a newly
# inserted duplicate of the "Cloud Curriculum" paper; row 25 is a newly appended record;
# every other row keeps its paper but gets an updated page number in column A).
$ws.Cells.Item(2, 1).Value = 10
$ws.Cells.Item(2, 2).Value = '2018_2c'
$ws.Cells.Item(2, 3).Value = '2c. The enemies within - inhibitors to learning'
$ws.Cells.Item(2, 4).Value = 'policy setting'

$ws.Cells.Item(3, 1).Value = 48
$ws.Cells.Item(3, 2).Value = '2019_3294016.3294021'
$ws.Cells.Item(3, 3).Value = 'Data Protection and Privacy Regulations as an  Inter-Active-Constructive Practice'
$ws.Cells.Item(3, 4).Value = 'cybersecurity'

$ws.Cells.Item(4, 1).Value = 59
$ws.Cells.Item(4, 2).Value = '2019_3294016.3294026'
$ws.Cells.Item(4, 3).Value = 'Engaging with computer science when solving tangible problems'
$ws.Cells.Item(4, 4).Value = 'engaging non-cs pupils in cs problems'

$ws.Cells.Item(5, 1).Value = 82
$ws.Cells.Item(5, 2).Value = '2020_3372356.3372362'
$ws.Cells.Item(5, 3).Value = 'Computer Science Degree Accreditation in the UK:  A Post-Shadbolt Review Update'
$ws.Cells.Item(5, 4).Value = 'professional body (BCS)'

$ws.Cells.Item(6, 1).Value = 84
$ws.Cells.Item(6, 2).Value = '2020_3372356.3372363'
$ws.Cells.Item(6, 3).Value = 'Building an Enhanced Student Experience Reflections from the Department of Computer Science at the University of Bath'
$ws.Cells.Item(6, 4).Value = 'student experience & work placement'

$ws.Cells.Item(7, 1).Value = 90
$ws.Cells.Item(7, 2).Value = '2020_3372356.3372366'
$ws.Cells.Item(7, 3).Value = 'Increasing academic diversity and inter-disciplinarity of  Computer Science in Higher Education'
$ws.Cells.Item(7, 4).Value = 'curriculum'

$ws.Cells.Item(8, 1).Value = 92
$ws.Cells.Item(8, 2).Value = '2020_3372356.3372367'
$ws.Cells.Item(8, 3).Value = 'Designing a Portfolio-Oriented Curriculum using Problem  Based Learning'
$ws.Cells.Item(8, 4).Value = 'curriculum'

$ws.Cells.Item(9, 1).Value = 100
$ws.Cells.Item(9, 2).Value = '2020_3372356.3372371'
$ws.Cells.Item(9, 3).Value = 'The PASS Effect  How Peer Assisted Study Sessions Contribute to a Computing Community'
$ws.Cells.Item(9, 4).Value = 'peer support'

$ws.Cells.Item(10, 1).Value = 108
$ws.Cells.Item(10, 2).Value = '2021_3437914.3437971'
$ws.Cells.Item(10, 3).Value = 'A Repository of Resources and Exemplars for the Cloud Curriculum'
$ws.Cells.Item(10, 4).Value = 'cloud computing'

$ws.Cells.Item(11, 1).Value = 109
$ws.Cells.Item(11, 2).Value = '2021_3437914.3437971'
$ws.Cells.Item(11, 3).Value = 'A Repository of Resources and Exemplars for the Cloud Curriculum'
$ws.Cells.Item(11, 4).Value = 'Cloud computing'

$ws.Cells.Item(12, 1).Value = 120
$ws.Cells.Item(12, 2).Value = '2021_3437914.3437977'
$ws.Cells.Item(12, 3).Value = 'Supporting Early-Career Academics in the UK Computer Science  Community'
$ws.Cells.Item(12, 4).Value = 'training for early year academics in teaching computing'

$ws.Cells.Item(13, 1).Value = 125
$ws.Cells.Item(13, 2).Value = '2022_3498343.3498344'
$ws.Cells.Item(13, 3).Value = 'Narrowing and Stretching: Addressing the Challenge of  Multi-track Programming'
$ws.Cells.Item(13, 4).Value = 'I dunno. It could be "programming", but I think that would be misleading'

$ws.Cells.Item(14, 1).Value = 134
$ws.Cells.Item(14, 2).Value = '2022_3498343.3498349'
$ws.Cells.Item(14, 3).Value = 'Co-Constructing a Community of Practice for Early-Career Computer Science Academics in the UK'
$ws.Cells.Item(14, 4).Value = 'training new CS lecturers'

$ws.Cells.Item(15, 1).Value = 135
$ws.Cells.Item(15, 2).Value = '2022_3498343.3498349'
$ws.Cells.Item(15, 3).Value = 'Co-Constructing a Community of Practice for Early-Career  Computer Science Academics in the UK'
$ws.Cells.Item(15, 4).Value = 'professional development'

$ws.Cells.Item(16, 1).Value = 137
$ws.Cells.Item(16, 2).Value = '2022_3498343.3498350'
$ws.Cells.Item(16, 3).Value = 'Application of AmazonWeb Services within teaching & learning at Coventry University Group'
$ws.Cells.Item(16, 4).Value = 'cloud computing '

$ws.Cells.Item(17, 1).Value = 138
$ws.Cells.Item(17, 2).Value = '2022_3498343.3498351'
$ws.Cells.Item(17, 3).Value = 'Assessing Knowledge and Skills in Foresnics with Alternative Assessment pathways'
$ws.Cells.Item(17, 4).Value = '"forensics" and/or "interdisciplinary"'

$ws.Cells.Item(18, 1).Value = 140
$ws.Cells.Item(18, 2).Value = '2022_3498343.3498353'
$ws.Cells.Item(18, 3).Value = 'LMC+ Scratch: A recipe to construct a mental model of program execution'
$ws.Cells.Item(18, 4).Value = 'This could be "programming" but a more refined classifcation would be "Mental Models"'

$ws.Cells.Item(19, 1).Value = 151
$ws.Cells.Item(19, 2).Value = '2023_3573260.3573265'
$ws.Cells.Item(19, 3).Value = 'Data Science Course Design for a Large-Scale Cohort using  Individual Project-Based Learning'
$ws.Cells.Item(19, 4).Value = 'data science'

$ws.Cells.Item(20, 1).Value = 163
$ws.Cells.Item(20, 2).Value = '2024_3633053.3633055'
$ws.Cells.Item(20, 3).Value = 'Institute of Coding in Wales Digital Skills Bootcamps – Micro-Credentials: A Pilot Project'
$ws.Cells.Item(20, 4).Value = 'upskilling'

$ws.Cells.Item(21, 1).Value = 191
$ws.Cells.Item(21, 2).Value = '2025_3702212.3702217'
$ws.Cells.Item(21, 3).Value = 'Themes in the Declared Use of Generative Artificial Intelligence  in Assessment'
$ws.Cells.Item(21, 4).Value = 'Forensics'

$ws.Cells.Item(22, 1).Value = 197
$ws.Cells.Item(22, 2).Value = '2025_3702212.3702220'
$ws.Cells.Item(22, 3).Value = 'Enhancing conceptual understanding in early years of Computing  education'
$ws.Cells.Item(22, 4).Value = 'conceptual understanding (pre=programming?)'

$ws.Cells.Item(23, 1).Value = 199
$ws.Cells.Item(23, 2).Value = '2025_3702212.3702222'
$ws.Cells.Item(23, 3).Value = 'Enhancing Learning and Teaching Experience for International  Students in Computing Subjects'
$ws.Cells.Item(23, 4).Value = 'international students'

$ws.Cells.Item(24, 1).Value = 203
$ws.Cells.Item(24, 2).Value = '2025_3702212.3702224'
$ws.Cells.Item(24, 3).Value = 'Where Have All the Papers Gone? Priming the pump of  pedagogical publishing in Europe'
$ws.Cells.Item(24, 4).Value = 'staff development'

$ws.Cells.Item(25, 1).Value = 205
$ws.Cells.Item(25, 2).Value = '2025_3702212.3702225'
$ws.Cells.Item(25, 3).Value = 'Integrating Socially Responsible Computing Competencies in Computer Science and Software Engineering Education'
$ws.Cells.Item(25, 4).Value = 'socially responsible computing (needs to be broader than "ethics", I think)'
